# Insert a new weekly price record as row 534, pushing the existing
# rows 534:551 down to 535:552 (dimension grows from A1:R551 to A1:R552).
# Rows.Insert() shifts the data down and carries the row-534 formatting
# (incl. the date style on column D) along with it, so both the new row
# and the rows pushed down keep the correct number formats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(534).Insert()

$ws.Cells.Item(534, 1).Value = 3
$ws.Cells.Item(534, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(534, 3).Value = "Coquimbo"
$ws.Cells.Item(534, 4).Value = 45075
$ws.Cells.Item(534, 5).Value = 5
$ws.Cells.Item(534, 6).Value = 100112009
$ws.Cells.Item(534, 7).Value = "Acelga"
$ws.Cells.Item(534, 8).Value = "Sin especificar"
$ws.Cells.Item(534, 9).Value = "Primera"
$ws.Cells.Item(534, 10).Value = 230
$ws.Cells.Item(534, 11).Value = 3300
$ws.Cells.Item(534, 12).Value = 3500
$ws.Cells.Item(534, 13).Value = 3396
$ws.Cells.Item(534, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(534, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(534, 16).Value = 566
$ws.Cells.Item(534, 17).Value = 6
$ws.Cells.Item(534, 18).Value = "Hortaliza"
